# Bump the test-user fixtures used by the LumaShop registration tests
# (user09/tu09@maildrop.cc, user08/tu08@maildrop.cc -> user10/tu10@maildrop.cc,
#  user11/tu11@maildrop.cc) and leave the cursor on the next free scratch
# cell, mirroring the manual edit that was made while adding the
# reporting/screenshot-embedding support for failed test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here: new shared-string entries are appended in the order
# cells are written, so write C3 / B2 / B3 / C2 in this exact sequence to
# reproduce the target shared string table layout.
$ws.Range("C3").Value = "tu11@maildrop.cc"
$ws.Range("B2").Value = "user10"
$ws.Range("B3").Value = "user11"
$ws.Range("C2").Value = "tu10@maildrop.cc"

# Move/leave the selection on I8, as recorded in the saved workbook.
$excel.Goto($ws.Range("I8"))
